$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9838
$ws.Range("F3").Value = 417
$ws.Range("F9").Value = 721
$ws.Range("F11").Value = 1215
$ws.Range("F13").Value = 3046
$ws.Range("F14").Value = 2302
$ws.Range("F16").Value = 1998
$ws.Range("F20").Value = 1570
$ws.Range("F21").Value = 325
$ws.Range("F22").Value = 36
$ws.Range("F23").Value = 204
$ws.Range("F28").Value = 334
$ws.Range("F30").Value = 41
$ws.Range("F31").Value = 182
$ws.Range("F33").Value = 239
$ws.Range("F34").Value = 1572
$ws.Range("F36").Value = 380
$ws.Range("F38").Value = 408
$ws.Range("F39").Value = 848
$ws.Range("F41").Value = 329

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 26

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9838
$ws.Range("F5").Value = 26
$ws.Range("F11").Value = 721
$ws.Range("F13").Value = 1215
$ws.Range("F15").Value = 3046
$ws.Range("F16").Value = 2302
$ws.Range("F18").Value = 1998
$ws.Range("F22").Value = 1570
$ws.Range("F23").Value = 325
$ws.Range("F24").Value = 36
$ws.Range("F25").Value = 204
$ws.Range("F30").Value = 334
$ws.Range("F35").Value = 41
$ws.Range("F36").Value = 182
$ws.Range("F39").Value = 239
$ws.Range("F40").Value = 1572
$ws.Range("F43").Value = 380
$ws.Range("F45").Value = 408
$ws.Range("F46").Value = 848
$ws.Range("F48").Value = 329
